$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q3" sheet by copying "2022-Q2" (keeps header/border
#     styles + column layout identical), then trim & overwrite with Q3 data. ---
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($wsQ2)
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

# Old sheet had 31 data rows (rows 2-32); new one only needs 22 (rows 2-23).
$wsQ3.Range("A24:H32").Clear()

# Columns D:G hold numeric-looking values that must stay TEXT (matches the rest of the
# workbook, where every fund-sheet data column besides A/H is stored as a string).
$wsQ3.Range("D2:G23").NumberFormat = "@"

# --- Step 2: overwrite cell-by-cell with 2022-Q3 figures ---
$wsQ3.Cells.Item(2, 1).Value = 0
$wsQ3.Cells.Item(2, 2).Value = "040015"
$wsQ3.Cells.Item(2, 3).Value = "华安动态灵活配置混合A"
$wsQ3.Cells.Item(2, 4).Value = "22.58"
$wsQ3.Cells.Item(2, 5).Value = "77.56"
$wsQ3.Cells.Item(2, 6).Value = "3.91"
$wsQ3.Cells.Item(2, 7).Value = "0.8829"
$wsQ3.Cells.Item(2, 8).Value = 6
$wsQ3.Cells.Item(3, 1).Value = 1
$wsQ3.Cells.Item(3, 2).Value = "010699"
$wsQ3.Cells.Item(3, 3).Value = "东方红创新趋势混合"
$wsQ3.Cells.Item(3, 4).Value = "26.13"
$wsQ3.Cells.Item(3, 5).Value = "71.44"
$wsQ3.Cells.Item(3, 6).Value = "2.79"
$wsQ3.Cells.Item(3, 7).Value = "0.7290"
$wsQ3.Cells.Item(3, 8).Value = 8
$wsQ3.Cells.Item(4, 1).Value = 2
$wsQ3.Cells.Item(4, 2).Value = "040001"
$wsQ3.Cells.Item(4, 3).Value = "华安创新混合"
$wsQ3.Cells.Item(4, 4).Value = "15.29"
$wsQ3.Cells.Item(4, 5).Value = "72.21"
$wsQ3.Cells.Item(4, 6).Value = "3.63"
$wsQ3.Cells.Item(4, 7).Value = "0.5550"
$wsQ3.Cells.Item(4, 8).Value = 6
$wsQ3.Cells.Item(5, 1).Value = 3
$wsQ3.Cells.Item(5, 2).Value = "014007"
$wsQ3.Cells.Item(5, 3).Value = "华安制造升级一年持有混合A"
$wsQ3.Cells.Item(5, 4).Value = "12.31"
$wsQ3.Cells.Item(5, 5).Value = "92.39"
$wsQ3.Cells.Item(5, 6).Value = "4.05"
$wsQ3.Cells.Item(5, 7).Value = "0.4986"
$wsQ3.Cells.Item(5, 8).Value = 9
$wsQ3.Cells.Item(6, 1).Value = 4
$wsQ3.Cells.Item(6, 2).Value = "008271"
$wsQ3.Cells.Item(6, 3).Value = "大成优势企业混合A"
$wsQ3.Cells.Item(6, 4).Value = "9.60"
$wsQ3.Cells.Item(6, 5).Value = "79.85"
$wsQ3.Cells.Item(6, 6).Value = "4.73"
$wsQ3.Cells.Item(6, 7).Value = "0.4541"
$wsQ3.Cells.Item(6, 8).Value = 9
$wsQ3.Cells.Item(7, 1).Value = 5
$wsQ3.Cells.Item(7, 2).Value = "010792"
$wsQ3.Cells.Item(7, 3).Value = "华安成长先锋混合A"
$wsQ3.Cells.Item(7, 4).Value = "11.21"
$wsQ3.Cells.Item(7, 5).Value = "91.59"
$wsQ3.Cells.Item(7, 6).Value = "4.02"
$wsQ3.Cells.Item(7, 7).Value = "0.4506"
$wsQ3.Cells.Item(7, 8).Value = 9
$wsQ3.Cells.Item(8, 1).Value = 6
$wsQ3.Cells.Item(8, 2).Value = "006154"
$wsQ3.Cells.Item(8, 3).Value = "华安制造先锋混合A"
$wsQ3.Cells.Item(8, 4).Value = "10.65"
$wsQ3.Cells.Item(8, 5).Value = "93.29"
$wsQ3.Cells.Item(8, 6).Value = "3.97"
$wsQ3.Cells.Item(8, 7).Value = "0.4228"
$wsQ3.Cells.Item(8, 8).Value = 10
$wsQ3.Cells.Item(9, 1).Value = 7
$wsQ3.Cells.Item(9, 2).Value = "001487"
$wsQ3.Cells.Item(9, 3).Value = "宝盈优势产业灵活配置混合A"
$wsQ3.Cells.Item(9, 4).Value = "10.11"
$wsQ3.Cells.Item(9, 5).Value = "91.85"
$wsQ3.Cells.Item(9, 6).Value = "2.88"
$wsQ3.Cells.Item(9, 7).Value = "0.2912"
$wsQ3.Cells.Item(9, 8).Value = 7
$wsQ3.Cells.Item(10, 1).Value = 8
$wsQ3.Cells.Item(10, 2).Value = "014389"
$wsQ3.Cells.Item(10, 3).Value = "华安产业动力6个月持有混合A"
$wsQ3.Cells.Item(10, 4).Value = "6.79"
$wsQ3.Cells.Item(10, 5).Value = "93.58"
$wsQ3.Cells.Item(10, 6).Value = "4.01"
$wsQ3.Cells.Item(10, 7).Value = "0.2723"
$wsQ3.Cells.Item(10, 8).Value = 8
$wsQ3.Cells.Item(11, 1).Value = 9
$wsQ3.Cells.Item(11, 2).Value = "013619"
$wsQ3.Cells.Item(11, 3).Value = "华安动态灵活配置混合C"
$wsQ3.Cells.Item(11, 4).Value = "4.47"
$wsQ3.Cells.Item(11, 5).Value = "77.56"
$wsQ3.Cells.Item(11, 6).Value = "3.91"
$wsQ3.Cells.Item(11, 7).Value = "0.1748"
$wsQ3.Cells.Item(11, 8).Value = 6
$wsQ3.Cells.Item(12, 1).Value = 10
$wsQ3.Cells.Item(12, 2).Value = "010793"
$wsQ3.Cells.Item(12, 3).Value = "华安成长先锋混合C"
$wsQ3.Cells.Item(12, 4).Value = "3.40"
$wsQ3.Cells.Item(12, 5).Value = "91.59"
$wsQ3.Cells.Item(12, 6).Value = "4.02"
$wsQ3.Cells.Item(12, 7).Value = "0.1367"
$wsQ3.Cells.Item(12, 8).Value = 9
$wsQ3.Cells.Item(13, 1).Value = 11
$wsQ3.Cells.Item(13, 2).Value = "012771"
$wsQ3.Cells.Item(13, 3).Value = "宝盈优势产业灵活配置混合C"
$wsQ3.Cells.Item(13, 4).Value = "3.62"
$wsQ3.Cells.Item(13, 5).Value = "91.85"
$wsQ3.Cells.Item(13, 6).Value = "2.88"
$wsQ3.Cells.Item(13, 7).Value = "0.1043"
$wsQ3.Cells.Item(13, 8).Value = 7
$wsQ3.Cells.Item(14, 1).Value = 12
$wsQ3.Cells.Item(14, 2).Value = "010738"
$wsQ3.Cells.Item(14, 3).Value = "大成优选升级一年持有期混合A"
$wsQ3.Cells.Item(14, 4).Value = "3.65"
$wsQ3.Cells.Item(14, 5).Value = "69.50"
$wsQ3.Cells.Item(14, 6).Value = "2.62"
$wsQ3.Cells.Item(14, 7).Value = "0.0956"
$wsQ3.Cells.Item(14, 8).Value = 10
$wsQ3.Cells.Item(15, 1).Value = 13
$wsQ3.Cells.Item(15, 2).Value = "002707"
$wsQ3.Cells.Item(15, 3).Value = "摩根士丹利华鑫科技领先灵活配置混合A"
$wsQ3.Cells.Item(15, 4).Value = "1.75"
$wsQ3.Cells.Item(15, 5).Value = "94.13"
$wsQ3.Cells.Item(15, 6).Value = "4.77"
$wsQ3.Cells.Item(15, 7).Value = "0.0835"
$wsQ3.Cells.Item(15, 8).Value = 4
$wsQ3.Cells.Item(16, 1).Value = 14
$wsQ3.Cells.Item(16, 2).Value = "008272"
$wsQ3.Cells.Item(16, 3).Value = "大成优势企业混合C"
$wsQ3.Cells.Item(16, 4).Value = "1.16"
$wsQ3.Cells.Item(16, 5).Value = "79.85"
$wsQ3.Cells.Item(16, 6).Value = "4.73"
$wsQ3.Cells.Item(16, 7).Value = "0.0549"
$wsQ3.Cells.Item(16, 8).Value = 9
$wsQ3.Cells.Item(17, 1).Value = 15
$wsQ3.Cells.Item(17, 2).Value = "014008"
$wsQ3.Cells.Item(17, 3).Value = "华安制造升级一年持有混合C"
$wsQ3.Cells.Item(17, 4).Value = "0.61"
$wsQ3.Cells.Item(17, 5).Value = "92.39"
$wsQ3.Cells.Item(17, 6).Value = "4.05"
$wsQ3.Cells.Item(17, 7).Value = "0.0247"
$wsQ3.Cells.Item(17, 8).Value = 9
$wsQ3.Cells.Item(18, 1).Value = 16
$wsQ3.Cells.Item(18, 2).Value = "014390"
$wsQ3.Cells.Item(18, 3).Value = "华安产业动力6个月持有混合C"
$wsQ3.Cells.Item(18, 4).Value = "0.59"
$wsQ3.Cells.Item(18, 5).Value = "93.58"
$wsQ3.Cells.Item(18, 6).Value = "4.01"
$wsQ3.Cells.Item(18, 7).Value = "0.0237"
$wsQ3.Cells.Item(18, 8).Value = 8
$wsQ3.Cells.Item(19, 1).Value = 17
$wsQ3.Cells.Item(19, 2).Value = "013507"
$wsQ3.Cells.Item(19, 3).Value = "华安制造先锋混合C"
$wsQ3.Cells.Item(19, 4).Value = "0.52"
$wsQ3.Cells.Item(19, 5).Value = "93.29"
$wsQ3.Cells.Item(19, 6).Value = "3.97"
$wsQ3.Cells.Item(19, 7).Value = "0.0206"
$wsQ3.Cells.Item(19, 8).Value = 10
$wsQ3.Cells.Item(20, 1).Value = 18
$wsQ3.Cells.Item(20, 2).Value = "007316"
$wsQ3.Cells.Item(20, 3).Value = "交银施罗德可转债债券A"
$wsQ3.Cells.Item(20, 4).Value = "0.66"
$wsQ3.Cells.Item(20, 5).Value = "23.69"
$wsQ3.Cells.Item(20, 6).Value = "0.76"
$wsQ3.Cells.Item(20, 7).Value = "0.0050"
$wsQ3.Cells.Item(20, 8).Value = 10
$wsQ3.Cells.Item(21, 1).Value = 19
$wsQ3.Cells.Item(21, 2).Value = "014871"
$wsQ3.Cells.Item(21, 3).Value = "摩根士丹利华鑫科技领先灵活配置混合C"
$wsQ3.Cells.Item(21, 4).Value = "0.08"
$wsQ3.Cells.Item(21, 5).Value = "94.13"
$wsQ3.Cells.Item(21, 6).Value = "4.77"
$wsQ3.Cells.Item(21, 7).Value = "0.0038"
$wsQ3.Cells.Item(21, 8).Value = 4
$wsQ3.Cells.Item(22, 1).Value = 20
$wsQ3.Cells.Item(22, 2).Value = "007317"
$wsQ3.Cells.Item(22, 3).Value = "交银施罗德可转债债券C"
$wsQ3.Cells.Item(22, 4).Value = "0.38"
$wsQ3.Cells.Item(22, 5).Value = "23.69"
$wsQ3.Cells.Item(22, 6).Value = "0.76"
$wsQ3.Cells.Item(22, 7).Value = "0.0029"
$wsQ3.Cells.Item(22, 8).Value = 10
$wsQ3.Cells.Item(23, 1).Value = 21
$wsQ3.Cells.Item(23, 2).Value = "010739"
$wsQ3.Cells.Item(23, 3).Value = "大成优选升级一年持有期混合C"
$wsQ3.Cells.Item(23, 4).Value = "0.08"
$wsQ3.Cells.Item(23, 5).Value = "69.50"
$wsQ3.Cells.Item(23, 6).Value = "2.62"
$wsQ3.Cells.Item(23, 7).Value = "0.0021"
$wsQ3.Cells.Item(23, 8).Value = 10
# --- Step 3: "总计" (summary) sheet - add the 2022-Q3 row at the top of the data
#     (row 2) and shift the existing quarters down by one row. ---
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 22
$wsTotal.Cells.Item(2, 4).Value = 5.29
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q2"
$wsTotal.Cells.Item(3, 3).Value = 31
$wsTotal.Cells.Item(3, 4).Value = 8.52
$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(4, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(4, 3).Value = 19
$wsTotal.Cells.Item(4, 4).Value = 6.72
$wsTotal.Cells.Item(5, 1).Value = 3
$wsTotal.Cells.Item(5, 2).Value = "2021-Q4"
$wsTotal.Cells.Item(5, 3).Value = 39
$wsTotal.Cells.Item(5, 4).Value = 10.75
$wsTotal.Cells.Item(6, 1).Value = 4
$wsTotal.Cells.Item(6, 2).Value = "2021-Q3"
$wsTotal.Cells.Item(6, 3).Value = 12
$wsTotal.Cells.Item(6, 4).Value = 5.98
$wsTotal.Cells.Item(7, 1).Value = 5
$wsTotal.Cells.Item(7, 2).Value = "2021-Q2"
$wsTotal.Cells.Item(7, 3).Value = 9
$wsTotal.Cells.Item(7, 4).Value = 3.41
$wsTotal.Cells.Item(8, 1).Value = 6
$wsTotal.Cells.Item(8, 2).Value = "2021-Q1"
$wsTotal.Cells.Item(8, 3).Value = 2
$wsTotal.Cells.Item(8, 4).Value = 0.73